# Ajout d'une nouvelle colonne de présence (CA) pour la date du 2025-11-15
# (numéro de série Excel 45976), juste après la dernière colonne existante (BZ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Nouvel en-tête de date.
$ws.Range("CA1").Value = 45976

# 2) Valeurs de présence de chaque joueur pour cette date.
$ws.Range("CA2").Value = "P"
$ws.Range("CA3").Value = "P"
$ws.Range("CA4").Value = "P"
$ws.Range("CA5").Value = "B"
$ws.Range("CA6").Value = "B"
$ws.Range("CA7").Value = "P"
$ws.Range("CA8").Value = "P"
$ws.Range("CA9").Value = "P"
$ws.Range("CA10").Value = "P"
$ws.Range("CA11").Value = "P"
# CA12 : le joueur de la ligne 12 est parti en cours de saison (la ligne
# s'arrête à la colonne AX) ; pas de nouvelle cellule pour lui.
$ws.Range("CA13").Value = "B"
$ws.Range("CA14").Value = "P"
$ws.Range("CA15").Value = "B"
$ws.Range("CA16").Value = "P"
$ws.Range("CA17").Value = "P"
$ws.Range("CA18").Value = "P"
$ws.Range("CA19").Value = "P"
$ws.Range("CA20").Value = "P"
# CA21 : joueur parti, cellule laissée vide (seul le style sera copié plus bas).
$ws.Range("CA22").Value = "P"
$ws.Range("CA23").Value = "B"
$ws.Range("CA24").Value = "P"
$ws.Range("CA25").Value = "P"
$ws.Range("CA26").Value = "P"
$ws.Range("CA27").Value = "P"
$ws.Range("CA28").Value = "P"
$ws.Range("CA29").Value = "RH"

# 3) Copier la mise en forme de la dernière colonne de présence (BZ) vers la
#    nouvelle colonne (CA) pour toutes les lignes du tableau, afin de garder
#    les mêmes styles (format de date centré pour l'en-tête, texte centré
#    pour les cellules de présence) sans perturber les valeurs déjà saisies.
$ws.Range("BZ1:BZ29").Copy()
$ws.Range("CA1:CA29").PasteSpecial(-4122)

# La ligne 12 ne doit pas avoir de cellule CA12 du tout.
$ws.Range("CA12").Clear()

# 4) Mettre à jour la cellule sélectionnée, comme dans le classeur d'origine.
$ws.Range("CC24").Select() | Out-Null
